# Update portfolio data for all three sheets (re-fetched 2025-09-21 08:55).
# Rows are re-sorted ascending by stock code, a new "现金" (cash) row is
# appended at the end of each sheet, and every row's modification timestamp
# moves from 202509211628 to 202509211655.

function Set-SmartCell($ws, $row, $col, $val) {
    # Mirrors how a human would type data into Excel: a cell whose text is
    # purely numeric (stock codes like "000089", timestamps like
    # "202509211655") needs to be forced to Text so leading zeros / the
    # full digit string survive; anything else (names, mixed alnum codes
    # like "HK01810") is left on General, which Excel already stores as
    # text on its own.
    $cell = $ws.Cells.Item($row, $col)
    if ($val -is [string] -and $val -match '^[0-9]+(\.[0-9]+)?$') {
        $cell.NumberFormat = "@"
    }
    $cell.Value = $val
}

function Write-PortfolioSheet($ws, $rows, $lastCol) {
    $lastRow = $rows.Length + 1
    # Wipe the old data block (header row untouched) before rewriting, so
    # stale rows beyond the new row count (shouldn't happen here, but keeps
    # this robust) don't linger.
    $clearRange = $ws.Range($ws.Cells.Item(2, 1), $ws.Cells.Item($lastRow + 5, $lastCol))
    $clearRange.ClearContents()

    for ($i = 0; $i -lt $rows.Length; $i++) {
        $row = $rows[$i]
        $r = $i + 2
        for ($c = 0; $c -lt $row.Length; $c++) {
            $col = $c + 1
            $val = $row[$c]
            if ($val -eq "") {
                # Blank source/suggested-ratio cells on sheet3 — leave empty.
                continue
            }
            Set-SmartCell $ws $r $col $val
        }
    }
}

# NOTE: array entries use a *trailing* comma at end-of-line. The leading-
# comma continuation style (comma at the start of the next line) is
# mis-parsed by this host's PowerShell engine -- it silently flattens the
# first row into loose top-level elements instead of keeping it nested --
# so every multi-row literal below intentionally avoids that form.
$sheet1Rows = @(
    @("大智 (稳健智远)", "000089", "深圳机场", 5.03, "202509211655"),
    @("大智 (稳健智远)", "000333", "美的集团", 9.71, "202509211655"),
    @("大智 (稳健智远)", "000831", "中国稀土", 9.21, "202509211655"),
    @("大智 (稳健智远)", "510300", "沪深300ETF", 5.26, "202509211655"),
    @("大智 (稳健智远)", "513400", "道琼斯ETF", 5.11, "202509211655"),
    @("大智 (稳健智远)", "518880", "黄金ETF", 2.06, "202509211655"),
    @("大智 (稳健智远)", "600085", "同仁堂", 1.92, "202509211655"),
    @("大智 (稳健智远)", "601899", "紫金矿业", 0.98, "202509211655"),
    @("大智 (稳健智远)", "100000", "现金", 60.73, "202509211655")
)

$sheet2Rows = @(
    @("大成 (锐进先锋)", "000725", "京东方A", 4.91, "202509211655"),
    @("大成 (锐进先锋)", "001380", "华纬科技", 5.22, "202509211655"),
    @("大成 (锐进先锋)", "002074", "国轩高科", 4.75, "202509211655"),
    @("大成 (锐进先锋)", "159781", "科创创业ETF", 6.11, "202509211655"),
    @("大成 (锐进先锋)", "513100", "纳指ETF", 5.17, "202509211655"),
    @("大成 (锐进先锋)", "513290", "纳指生物科技ETF", 0.97, "202509211655"),
    @("大成 (锐进先锋)", "600580", "卧龙电驱", 5.69, "202509211655"),
    @("大成 (锐进先锋)", "601878", "浙商证券", 4.89, "202509211655"),
    @("大成 (锐进先锋)", "603119", "浙江荣泰", 0.03, "202509211655"),
    @("大成 (锐进先锋)", "HK01810", "小米集团-W", 1.02, "202509211655"),
    @("大成 (锐进先锋)", "100000", "现金", 61.24, "202509211655")
)

$sheet3Rows = @(
    @("范式进化投资组合", "000089", "深圳机场", "", "", 5.04, "202509211655"),
    @("范式进化投资组合", "000333", "美的集团", "", "", 1.01, "202509211655"),
    @("范式进化投资组合", "000725", "京东方A", "", "", 5.05, "202509211655"),
    @("范式进化投资组合", "000831", "中国稀土", "", "", 9.6, "202509211655"),
    @("范式进化投资组合", "159781", "科创创业ETF", "", "", 6.85, "202509211655"),
    @("范式进化投资组合", "510050", "上证50ETF", "", "", 1.01, "202509211655"),
    @("范式进化投资组合", "510300", "沪深300ETF", "", "", 5.38, "202509211655"),
    @("范式进化投资组合", "513100", "纳指ETF", "", "", 3.14, "202509211655"),
    @("范式进化投资组合", "513290", "纳指生物科技ETF", "", "", 0.98, "202509211655"),
    @("范式进化投资组合", "513400", "道琼斯ETF", "", "", 5.05, "202509211655"),
    @("范式进化投资组合", "518880", "黄金ETF", "", "", 1, "202509211655"),
    @("范式进化投资组合", "600085", "同仁堂", "", "", 0.98, "202509211655"),
    @("范式进化投资组合", "100000", "现金", "", "", 54.9, "202509211655")
)

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
Write-PortfolioSheet $ws1 $sheet1Rows 5

$ws2 = $wb.Worksheets.Item(2)
Write-PortfolioSheet $ws2 $sheet2Rows 5

$ws3 = $wb.Worksheets.Item(3)
Write-PortfolioSheet $ws3 $sheet3Rows 7
